# "Add guardian tab additional column"
#
# The Guardian worksheet gains a new column H ("DateOfBirth") containing a
# sample date of birth value for the existing guardian row, matching the
# style already used for dates elsewhere in the workbook (Patient!D2).

$wb = $excel.ActiveWorkbook

$guardian = $wb.Worksheets.Item("Guardian")
$patient  = $wb.Worksheets.Item("Patient")

# Make sure Patient's selection (D2) is recorded before we switch away from
# it, then finish on the Guardian tab so it ends up as the active sheet.
[void]$patient.Select()
[void]$patient.Range("D2").Select()

[void]$guardian.Select()

# New header cell for the added column.
$guardian.Range("H1").Value = "DateOfBirth"

# New data cell: a date of birth value (1/22/1974 -> serial 27051).
$guardian.Range("H2").Value = 27051

# Reuse the date number format already used on the Patient sheet (D2) so we
# pick up the existing style (numFmtId 14, "m/d/yyyy") instead of Excel
# fabricating a brand-new style entry.
[void]$patient.Range("D2").Copy()
[void]$guardian.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the selection on the newly added cell, matching the saved view state.
[void]$guardian.Range("H2").Select()
